# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll        = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet (sheet1) ---
$wsExhibition.Range("F2").Value = 121
$wsExhibition.Range("G3").Value = "已售罄"
$wsExhibition.Range("F4").Value = 213
$wsExhibition.Range("F10").Value = 36
$wsExhibition.Range("F11").Value = 6986
$wsExhibition.Range("F12").Value = 250
$wsExhibition.Range("F13").Value = 396
$wsExhibition.Range("F14").Value = 3409
$wsExhibition.Range("F15").Value = 241
$wsExhibition.Range("F16").Value = 435
$wsExhibition.Range("F17").Value = 263
$wsExhibition.Range("F18").Value = 579
$wsExhibition.Range("F19").Value = 54

# --- 全部类型 sheet (sheet4) ---
$wsAll.Range("F4").Value = 121
$wsAll.Range("G5").Value = "已售罄"
$wsAll.Range("F6").Value = 213
$wsAll.Range("F12").Value = 36
$wsAll.Range("F14").Value = 6986
$wsAll.Range("F16").Value = 250
$wsAll.Range("F17").Value = 396
$wsAll.Range("F18").Value = 3409
$wsAll.Range("F19").Value = 241
$wsAll.Range("F20").Value = 435
$wsAll.Range("F21").Value = 263
$wsAll.Range("F22").Value = 579
$wsAll.Range("F23").Value = 54
